$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 166.35294
$ws.Cells.Item(2, 9).Value = 159.86667
$ws.Cells.Item(2, 10).Value = 215
$ws.Cells.Item(2, 11).Value = 159.86667
$ws.Cells.Item(2, 12).Value = 215
$ws.Cells.Item(2, 13).Value = -46.86667
$ws.Cells.Item(2, 14).Value = -441

$ws.Cells.Item(15, 8).Value = 798.14
$ws.Cells.Item(15, 9).Value = 798.14
$ws.Cells.Item(15, 11).Value = 2394.42
$ws.Cells.Item(15, 13).Value = -2225.42

$ws.Cells.Item(40, 8).Value = 2084.6
$ws.Cells.Item(40, 9).Value = 1844.6666
$ws.Cells.Item(40, 10).Value = 2280.9092
$ws.Cells.Item(40, 11).Value = 1844.6666
$ws.Cells.Item(40, 12).Value = 2280.9092
$ws.Cells.Item(40, 13).Value = -1669.6666
$ws.Cells.Item(40, 14).Value = -2630.9092

$ws.Cells.Item(44, 8).Value = 20000
$ws.Cells.Item(44, 10).Value = 20000
$ws.Cells.Item(44, 12).Value = 20000
$ws.Cells.Item(44, 14).Value = -20924

$ws.Cells.Item(98, 8).Value = 3470.739
$ws.Cells.Item(98, 9).Value = 3573.9546
$ws.Cells.Item(98, 10).Value = 1200
$ws.Cells.Item(98, 11).Value = 3573.9546
$ws.Cells.Item(98, 12).Value = 1200
$ws.Cells.Item(98, 13).Value = -2075.9546
$ws.Cells.Item(98, 14).Value = -4196

$ws.Cells.Item(112, 8).Value = 2333.087
$ws.Cells.Item(112, 9).Value = 619.6
$ws.Cells.Item(112, 10).Value = 2809.0557
$ws.Cells.Item(112, 11).Value = 1858.8
$ws.Cells.Item(112, 12).Value = 8427.167099999999
$ws.Cells.Item(112, 13).Value = -750.8000000000002
$ws.Cells.Item(112, 14).Value = -10643.1671

$ws.Cells.Item(122, 8).Value = 3470.739
$ws.Cells.Item(122, 9).Value = 3573.9546
$ws.Cells.Item(122, 10).Value = 1200
$ws.Cells.Item(122, 11).Value = 10721.8638
$ws.Cells.Item(122, 12).Value = 3600
$ws.Cells.Item(122, 13).Value = -8271.863799999999
$ws.Cells.Item(122, 14).Value = -8500

$ws.Cells.Item(125, 8).Value = 3167
$ws.Cells.Item(125, 9).Value = 2272.5
$ws.Cells.Item(125, 10).Value = 3933.7144
$ws.Cells.Item(125, 11).Value = 20452.5
$ws.Cells.Item(125, 12).Value = 35403.4296
$ws.Cells.Item(125, 13).Value = -17992.5
$ws.Cells.Item(125, 14).Value = -40323.4296

$ws.Cells.Item(138, 8).Value = 445327.8
$ws.Cells.Item(138, 10).Value = 589360.2
$ws.Cells.Item(138, 12).Value = 1768080.6
$ws.Cells.Item(138, 14).Value = -1778360.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5532.104
$ws.Cells.Item(32, 9).Value = 4880.427
$ws.Cells.Item(32, 11).Value = 4880.427
$ws.Cells.Item(32, 13).Value = -4593.427

$ws.Cells.Item(74, 8).Value = 3368.6667
$ws.Cells.Item(74, 9).Value = 2842.4
$ws.Cells.Item(74, 10).Value = 6000
$ws.Cells.Item(74, 11).Value = 2842.4
$ws.Cells.Item(74, 12).Value = 6000
$ws.Cells.Item(74, 13).Value = -1968.4
$ws.Cells.Item(74, 14).Value = -7748

$ws.Cells.Item(77, 8).Value = 3368.6667
$ws.Cells.Item(77, 9).Value = 2842.4
$ws.Cells.Item(77, 10).Value = 6000
$ws.Cells.Item(77, 11).Value = 14212
$ws.Cells.Item(77, 12).Value = 30000
$ws.Cells.Item(77, 13).Value = -9844
$ws.Cells.Item(77, 14).Value = -38736

$ws.Cells.Item(110, 8).Value = 1514.4615
$ws.Cells.Item(110, 9).Value = 906.125
$ws.Cells.Item(110, 10).Value = 2487.8
$ws.Cells.Item(110, 11).Value = 906.125
$ws.Cells.Item(110, 12).Value = 2487.8
$ws.Cells.Item(110, 13).Value = 1138.875
$ws.Cells.Item(110, 14).Value = -6577.8

$ws.Cells.Item(132, 8).Value = 2472.9546
$ws.Cells.Item(132, 9).Value = 2412.15
$ws.Cells.Item(132, 10).Value = 2523.625
$ws.Cells.Item(132, 11).Value = 7236.450000000001
$ws.Cells.Item(132, 12).Value = 7570.875
$ws.Cells.Item(132, 13).Value = -4706.450000000001
$ws.Cells.Item(132, 14).Value = -12630.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3105.3333
$ws.Cells.Item(20, 9).Value = 2898.4614
$ws.Cells.Item(20, 10).Value = 4450
$ws.Cells.Item(20, 11).Value = 2898.4614
$ws.Cells.Item(20, 12).Value = 4450
$ws.Cells.Item(20, 13).Value = -2651.4614
$ws.Cells.Item(20, 14).Value = -4944

$ws.Cells.Item(105, 8).Value = 111112890
$ws.Cells.Item(105, 9).Value = 125001736
$ws.Cells.Item(105, 11).Value = 125001736
$ws.Cells.Item(105, 13).Value = -124999989

$ws.Cells.Item(134, 8).Value = 1149.7646
$ws.Cells.Item(134, 9).Value = 969.8
$ws.Cells.Item(134, 11).Value = 2909.4
$ws.Cells.Item(134, 13).Value = -374.3999999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1756.3684
$ws.Cells.Item(31, 10).Value = 4999
$ws.Cells.Item(31, 12).Value = 4999
$ws.Cells.Item(31, 14).Value = -5589

$ws.Cells.Item(34, 8).Value = 1756.3684
$ws.Cells.Item(34, 10).Value = 4999
$ws.Cells.Item(34, 12).Value = 4999
$ws.Cells.Item(34, 14).Value = -5403

$ws.Cells.Item(62, 8).Value = 11766903
$ws.Cells.Item(62, 9).Value = 2362.5
$ws.Cells.Item(62, 10).Value = 40001800
$ws.Cells.Item(62, 11).Value = 2362.5
$ws.Cells.Item(62, 12).Value = 40001800
$ws.Cells.Item(62, 13).Value = -1738.5
$ws.Cells.Item(62, 14).Value = -40003048

$ws.Cells.Item(65, 8).Value = 11766903
$ws.Cells.Item(65, 9).Value = 2362.5
$ws.Cells.Item(65, 10).Value = 40001800
$ws.Cells.Item(65, 11).Value = 11812.5
$ws.Cells.Item(65, 12).Value = 200009000
$ws.Cells.Item(65, 13).Value = -8692.5
$ws.Cells.Item(65, 14).Value = -200015240

$ws.Cells.Item(94, 8).Value = 2042.3077
$ws.Cells.Item(94, 9).Value = 1799.5
$ws.Cells.Item(94, 10).Value = 2150.2222
$ws.Cells.Item(94, 11).Value = 1799.5
$ws.Cells.Item(94, 12).Value = 2150.2222
$ws.Cells.Item(94, 13).Value = -1348.5
$ws.Cells.Item(94, 14).Value = -3052.2222

$ws.Cells.Item(99, 8).Value = 1278.7273
$ws.Cells.Item(99, 9).Value = 1324
$ws.Cells.Item(99, 10).Value = 1158
$ws.Cells.Item(99, 11).Value = 1324
$ws.Cells.Item(99, 12).Value = 1158
$ws.Cells.Item(99, 13).Value = 174
$ws.Cells.Item(99, 14).Value = -4154

$ws.Cells.Item(105, 8).Value = 936.8333
$ws.Cells.Item(105, 9).Value = 922
$ws.Cells.Item(105, 11).Value = 922
$ws.Cells.Item(105, 13).Value = 825

$ws.Cells.Item(119, 8).Value = 22000
$ws.Cells.Item(119, 10).Value = 22000
$ws.Cells.Item(119, 12).Value = 22000
$ws.Cells.Item(119, 14).Value = -31676

$ws.Cells.Item(126, 8).Value = 1278.7273
$ws.Cells.Item(126, 9).Value = 1324
$ws.Cells.Item(126, 10).Value = 1158
$ws.Cells.Item(126, 11).Value = 3972
$ws.Cells.Item(126, 12).Value = 3474
$ws.Cells.Item(126, 13).Value = -1502
$ws.Cells.Item(126, 14).Value = -8414

$ws.Cells.Item(135, 8).Value = 32931.25
$ws.Cells.Item(135, 10).Value = 32931.25
$ws.Cells.Item(135, 12).Value = 32931.25
$ws.Cells.Item(135, 14).Value = -43071.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 555
$ws.Cells.Item(2, 9).Value = 26.461538
$ws.Cells.Item(2, 10).Value = 1536.5714
$ws.Cells.Item(2, 11).Value = 158.769228
$ws.Cells.Item(2, 12).Value = 9219.428400000001
$ws.Cells.Item(2, 13).Value = -45.769228
$ws.Cells.Item(2, 14).Value = -9445.428400000001

$ws.Cells.Item(8, 8).Value = 117.833336
$ws.Cells.Item(8, 9).Value = 117.833336
$ws.Cells.Item(8, 11).Value = 353.500008
$ws.Cells.Item(8, 13).Value = -214.500008

$ws.Cells.Item(32, 8).Value = 1750
$ws.Cells.Item(32, 9).Value = 2000
$ws.Cells.Item(32, 10).Value = 1500
$ws.Cells.Item(32, 11).Value = 6000
$ws.Cells.Item(32, 12).Value = 4500
$ws.Cells.Item(32, 13).Value = -5717
$ws.Cells.Item(32, 14).Value = -5066

$ws.Cells.Item(92, 8).Value = 237.02563
$ws.Cells.Item(92, 9).Value = 225.53334
$ws.Cells.Item(92, 10).Value = 275.33334
$ws.Cells.Item(92, 11).Value = 676.6000200000001
$ws.Cells.Item(92, 12).Value = 826.0000200000001
$ws.Cells.Item(92, 13).Value = 571.3999799999999
$ws.Cells.Item(92, 14).Value = -3322.00002

$ws.Cells.Item(93, 8).Value = 6755.5
$ws.Cells.Item(93, 10).Value = 6755.5
$ws.Cells.Item(93, 12).Value = 20266.5
$ws.Cells.Item(93, 14).Value = -24010.5

$ws.Cells.Item(94, 8).Value = 4224.385
$ws.Cells.Item(94, 10).Value = 4157.75
$ws.Cells.Item(94, 12).Value = 12473.25
$ws.Cells.Item(94, 14).Value = -13825.25

$ws.Cells.Item(96, 8).Value = 7757.143
$ws.Cells.Item(96, 10).Value = 7757.143
$ws.Cells.Item(96, 12).Value = 23271.429
$ws.Cells.Item(96, 14).Value = -27389.429

$ws.Cells.Item(105, 8).Value = 7423.077
$ws.Cells.Item(105, 10).Value = 7423.077
$ws.Cells.Item(105, 12).Value = 22269.231
$ws.Cells.Item(105, 14).Value = -27511.231

$ws.Cells.Item(131, 8).Value = 19231760
$ws.Cells.Item(131, 10).Value = 1081.9318
$ws.Cells.Item(131, 12).Value = 3245.7954
$ws.Cells.Item(131, 14).Value = -13325.7954

$ws.Cells.Item(139, 8).Value = 2201.3667
$ws.Cells.Item(139, 9).Value = 2729.077
$ws.Cells.Item(139, 11).Value = 8187.231000000001
$ws.Cells.Item(139, 13).Value = -3047.231000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 20000
$ws.Cells.Item(15, 10).Value = 20000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 14).Value = -20576

$ws.Cells.Item(81, 8).Value = 20000
$ws.Cells.Item(81, 10).Value = 20000
$ws.Cells.Item(81, 12).Value = 20000
$ws.Cells.Item(81, 14).Value = -21996

$ws.Cells.Item(84, 8).Value = 20000
$ws.Cells.Item(84, 10).Value = 20000
$ws.Cells.Item(84, 12).Value = 60000
$ws.Cells.Item(84, 14).Value = -69984

$ws.Cells.Item(122, 8).Value = 4185.353
$ws.Cells.Item(122, 9).Value = 4841.4
$ws.Cells.Item(122, 10).Value = 3248.1428
$ws.Cells.Item(122, 11).Value = 14524.2
$ws.Cells.Item(122, 12).Value = 9744.428400000001
$ws.Cells.Item(122, 13).Value = -12074.2
$ws.Cells.Item(122, 14).Value = -14644.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).ClearContents()  # M25 removed
$ws.Cells.Item(25, 14).ClearContents()  # N25 removed

$ws.Cells.Item(68, 8).Value = 2989.524
$ws.Cells.Item(68, 9).Value = 2998.9473
$ws.Cells.Item(68, 10).Value = 2900
$ws.Cells.Item(68, 11).Value = 2998.9473
$ws.Cells.Item(68, 12).Value = 2900
$ws.Cells.Item(68, 13).Value = -2249.9473
$ws.Cells.Item(68, 14).Value = -4398

$ws.Cells.Item(71, 8).Value = 2989.524
$ws.Cells.Item(71, 9).Value = 2998.9473
$ws.Cells.Item(71, 10).Value = 2900
$ws.Cells.Item(71, 11).Value = 14994.7365
$ws.Cells.Item(71, 12).Value = 14500
$ws.Cells.Item(71, 13).Value = -11250.7365
$ws.Cells.Item(71, 14).Value = -21988

$ws.Cells.Item(121, 8).Value = 30000
$ws.Cells.Item(121, 10).Value = 30000
$ws.Cells.Item(121, 12).Value = 30000
$ws.Cells.Item(121, 14).Value = -33494

$ws.Cells.Item(122, 8).Value = 50002400
$ws.Cells.Item(122, 9).Value = 62502250
$ws.Cells.Item(122, 11).Value = 187506750
$ws.Cells.Item(122, 13).Value = -187504300

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 339.4762
$ws.Cells.Item(113, 9).Value = 285.75
$ws.Cells.Item(113, 10).Value = 411.1111
$ws.Cells.Item(113, 11).Value = 857.25
$ws.Cells.Item(113, 12).Value = 1233.3333
$ws.Cells.Item(113, 13).Value = 1312.75
$ws.Cells.Item(113, 14).Value = -5573.3333

$ws.Cells.Item(121, 8).Value = 32500
$ws.Cells.Item(121, 10).Value = 32500
$ws.Cells.Item(121, 12).Value = 32500
$ws.Cells.Item(121, 14).Value = -35994

$ws.Cells.Item(132, 8).Value = 1858.4286
$ws.Cells.Item(132, 9).Value = 1672.7693
$ws.Cells.Item(132, 10).Value = 2394.7778
$ws.Cells.Item(132, 11).Value = 5018.3079
$ws.Cells.Item(132, 12).Value = 7184.3334
$ws.Cells.Item(132, 13).Value = -2488.3079
$ws.Cells.Item(132, 14).Value = -12244.3334
